# Edit script: add a new survey date column (16. 3. 2021) to both sheets
# and bump the "aktualizace" date in the footer notes from 9. 3. 2021 to 23. 3. 2021.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("data")
$ws2 = $wb.Worksheets.Item("pocetR")

# --- Sheet "data": add column Z (26) with header "16. 3. 2021" ---
# Copy formatting (border/font/alignment) from the previous header cell Y1 so the
# new header cell matches the existing header style.
$ws1.Range("Y1").Copy()
$ws1.Range("Z1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws1.Range("Z1").Value = "16. 3. 2021"

$ws1.Cells.Item(2,26).Value = 0.45
$ws1.Cells.Item(3,26).Value = 0.34
$ws1.Cells.Item(4,26).Value = 0.21
$ws1.Cells.Item(5,26).Value = 0.31
$ws1.Cells.Item(6,26).Value = 0.31
$ws1.Cells.Item(7,26).Value = 0.38
$ws1.Cells.Item(8,26).Value = 0.48
$ws1.Cells.Item(9,26).Value = 0.35
$ws1.Cells.Item(10,26).Value = 0.17
$ws1.Cells.Item(11,26).Value = 0.47
$ws1.Cells.Item(12,26).Value = 0.35
$ws1.Cells.Item(13,26).Value = 0.18
$ws1.Cells.Item(14,26).Value = 0.34
$ws1.Cells.Item(15,26).Value = 0.39
$ws1.Cells.Item(16,26).Value = 0.27
$ws1.Cells.Item(17,26).Value = 0.48
$ws1.Cells.Item(18,26).Value = 0.34
$ws1.Cells.Item(19,26).Value = 0.18
$ws1.Cells.Item(20,26).Value = 0.36
$ws1.Cells.Item(21,26).Value = 0.38
$ws1.Cells.Item(22,26).Value = 0.26
$ws1.Cells.Item(23,26).Value = 0.38
$ws1.Cells.Item(24,26).Value = 0.36
$ws1.Cells.Item(25,26).Value = 0.26
$ws1.Cells.Item(26,26).Value = 0.39
$ws1.Cells.Item(27,26).Value = 0.39
$ws1.Cells.Item(28,26).Value = 0.22
$ws1.Cells.Item(29,26).Value = 0.45
$ws1.Cells.Item(30,26).Value = 0.34
$ws1.Cells.Item(31,26).Value = 0.21
$ws1.Cells.Item(32,26).Value = 0.56
$ws1.Cells.Item(33,26).Value = 0.28
$ws1.Cells.Item(34,26).Value = 0.16
$ws1.Cells.Item(35,26).Value = 0.36
$ws1.Cells.Item(36,26).Value = 0.41
$ws1.Cells.Item(37,26).Value = 0.23
$ws1.Cells.Item(38,26).Value = 0.38
$ws1.Cells.Item(39,26).Value = 0.38
$ws1.Cells.Item(40,26).Value = 0.24
$ws1.Cells.Item(41,26).Value = 0.54
$ws1.Cells.Item(42,26).Value = 0.29
$ws1.Cells.Item(43,26).Value = 0.17
$ws1.Cells.Item(44,26).Value = 0.52
$ws1.Cells.Item(45,26).Value = 0.3
$ws1.Cells.Item(46,26).Value = 0.18
$ws1.Cells.Item(47,26).Value = 0.39
$ws1.Cells.Item(48,26).Value = 0.34
$ws1.Cells.Item(49,26).Value = 0.27
$ws1.Cells.Item(50,26).Value = 0.59
$ws1.Cells.Item(51,26).Value = 0.31
$ws1.Cells.Item(52,26).Value = 0.1
$ws1.Cells.Item(53,26).Value = 0.41
$ws1.Cells.Item(54,26).Value = 0.37
$ws1.Cells.Item(55,26).Value = 0.22
$ws1.Cells.Item(56,26).Value = 0.61
$ws1.Cells.Item(57,26).Value = 0.31
$ws1.Cells.Item(58,26).Value = 0.08
$ws1.Cells.Item(59,26).Value = 0.55
$ws1.Cells.Item(60,26).Value = 0.32
$ws1.Cells.Item(61,26).Value = 0.13

# --- Sheet "pocetR": add column Y (25) with header "16. 3. 2021" ---
$ws2.Range("X1").Copy()
$ws2.Range("Y1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws2.Range("Y1").Value = "16. 3. 2021"

$ws2.Cells.Item(2,25).Value = 1160
$ws2.Cells.Item(3,25).Value = 178
$ws2.Cells.Item(4,25).Value = 982
$ws2.Cells.Item(5,25).Value = 923
$ws2.Cells.Item(6,25).Value = 166
$ws2.Cells.Item(7,25).Value = 8
$ws2.Cells.Item(8,25).Value = 63
$ws2.Cells.Item(9,25).Value = 886
$ws2.Cells.Item(10,25).Value = 153
$ws2.Cells.Item(11,25).Value = 69
$ws2.Cells.Item(12,25).Value = 51
$ws2.Cells.Item(13,25).Value = 417
$ws2.Cells.Item(14,25).Value = 458
$ws2.Cells.Item(15,25).Value = 285
$ws2.Cells.Item(16,25).Value = 125
$ws2.Cells.Item(17,25).Value = 334
$ws2.Cells.Item(18,25).Value = 396
$ws2.Cells.Item(19,25).Value = 187
$ws2.Cells.Item(20,25).Value = 320
$ws2.Cells.Item(21,25).Value = 96
$ws2.Cells.Item(22,25).Value = 307
$ws2.Cells.Item(23,25).Value = 163
$ws2.Cells.Item(24,25).Value = 105

# Row 25 on "pocetR" is the trailing footer row; every other column in that row
# holds an empty string. Copy that same "empty string" cell state from the
# previous column (X25) into the new Y25 cell so it matches the existing pattern.
$ws2.Range("X25").Copy()
$ws2.Range("Y25").PasteSpecial(-4163)
$excel.CutCopyMode = $false

# --- Update the "aktualizace" date in both footer notes ---
$ws1.Range("A62").Value = "Život během pandemie, Obavy ze ztráty práce, % respondentů celkově a ve skupinách, aktualizace 23. 3. 2021"
$ws2.Range("A25").Value = "Život během pandemie, Obavy ze ztráty práce, velikost dotázaného souboru celkově a ve skupinách, aktualizace 23. 3. 2021"
